$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 382 (shifts rows 382:500 down to 383:501)
$ws.Rows(382).Insert()

# Populate the new row 382 with the new data record.
# Columns A,B,C,E,F,G,H,I,J,K,Q,T are constant across every record in this sheet.
$ws.Range("A382").Value = 5
$ws.Range("B382").Value = "Macroferia Regional de Talca"
$ws.Range("C382").Value = "Maule"
$ws.Range("D382").Value = 44627
$ws.Range("D382").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E382").Value = 7
$ws.Range("F382").Value = "Fruta"
$ws.Range("G382").Value = 100101
$ws.Range("H382").Value = "Berries"
$ws.Range("I382").Value = 100112025
$ws.Range("J382").Value = "Frutilla"
$ws.Range("K382").Value = "Sin especificar"
$ws.Range("L382").Value = "Primera"
$ws.Range("M382").Value = 80
$ws.Range("N382").Value = 6000
$ws.Range("O382").Value = 6000
$ws.Range("P382").Value = 6000
$ws.Range("Q382").Value = "$/bandeja 7 kilos"
$ws.Range("R382").Value = "Región del Maule"
$ws.Range("S382").Value = 857
$ws.Range("T382").Value = 7
